$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New values for columns B, C, D, E, G (rows 2-25).
# Column F ("Win") is left unchanged per the source diff.
$data = @(
    @(2, 0.6606524410359556, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 2.960089034096801),
    @(3, 0.1190320826869504, 0.306821227259698, 0.1494219747398047, 0.4942365360607697, 1.069511820747223),
    @(4, 1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 4.358119930609447),
    @(5, 0.6606524410359556, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 6.348428708163715),
    @(6, 1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 4.358119930609447),
    @(7, 0.003208871385164791, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 2.30264546444601),
    @(8, 3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694),
    @(9, 0.6606524410359556, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 2.960089034096801),
    @(10, 3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694),
    @(11, 0.1190320826869504, 0.306821227259698, 0.7527432677738641, 0.4942365360607697, 1.672833113781282),
    @(12, 3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694),
    @(13, 3.286832544864788, 1.655778082260271, 22.3905356188092, 0.4942365360607697, 27.82738278199502),
    @(14, 1.455362044514542, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 3.754798637575387),
    @(15, 0.6606524410359556, 0.306821227259698, 0.7527432677738641, 0.4942365360607697, 2.214453472130288),
    @(16, 3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 5.586269137925634),
    @(17, 3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694),
    @(18, 0.1190320826869504, 0.306821227259698, 0.1494219747398047, 0.4942365360607697, 1.069511820747223),
    @(19, 0.6606524410359556, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 6.348428708163715),
    @(20, 1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 4.358119930609447),
    @(21, 0.2917716402565462, 0.306821227259698, 0.7527432677738641, 0.4942365360607697, 1.845572671350878),
    @(22, 1.455362044514542, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 3.754798637575387),
    @(23, 1.455362044514542, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 3.754798637575387),
    @(24, 3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694),
    @(25, 3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
)

foreach ($r in $data) {
    $row = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]   # B
    $ws.Cells.Item($row, 3).Value = $r[2]   # C
    $ws.Cells.Item($row, 4).Value = $r[3]   # D
    $ws.Cells.Item($row, 5).Value = $r[4]   # E
    $ws.Cells.Item($row, 7).Value = $r[5]   # G
}

